# Update "想去人数" (want-to-go count) values in the F column on the
# "展览" and "全部类型" sheets to match the freshly generated data,
# as published by the gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览"     = @{ "F5" = 178; "F6" = 9518; "F10" = 1683; "F11" = 154; "F14" = 269; "F15" = 441; "F18" = 1309 }
    "全部类型" = @{ "F6" = 178; "F7" = 9518; "F11" = 1683; "F12" = 154; "F15" = 269; "F16" = 441; "F19" = 1309 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($cellRef in $cellUpdates.Keys) {
        $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
    }
}
